$d = $word.ActiveDocument
$found = $d.Content
$found.Find.Text = "N°"
$found.Find.Execute() | Out-Null
Write-Output ("start=" + $found.Start + " end=" + $found.End)
$xmlFrag = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
<w:proofErr w:type="spellStart"/>
<w:r><w:rPr><w:rFonts w:ascii="Arial Rounded MT Bold" w:hAnsi="Arial Rounded MT Bold"/></w:rPr><w:t>N°</w:t></w:r>
<w:proofErr w:type="spellEnd"/>
</w:p>
</w:body>
</w:document>
'@
try {
  $r = $found.InsertXML($xmlFrag)
  Write-Output ("insert ok " + $r)
} catch {
  Write-Output ("error: " + $_.Exception.Message)
}
